$wb = $excel.ActiveWorkbook

# "optimization_parameters" sheet (7th tab): remove the stray "Sheet" row
# (row 16 - a leftover label "Sheet" with values 3 / 4). Deleting it also
# drops the now-unused shared string "Sheet" and shifts the
# "simulation_timepoints" row up from row 17 to row 16.
$wsParams = $wb.Worksheets.Item(7)
$wsParams.Rows.Item(16).Delete() | Out-Null

# Update the lingering cell selections left on a few sheets.
$wsDegradation = $wb.Worksheets.Item(2)
$wsDegradation.Range("D40").Select() | Out-Null

$wsNetworkWeights = $wb.Worksheets.Item(6)
$wsNetworkWeights.Range("C7").Select() | Out-Null

$wsParams.Range("A16:XFD16").Select() | Out-Null

# The active/selected tab moves from "optimization_parameters" to the last
# sheet, "optimization_diagnostics".
$wsDiagnostics = $wb.Worksheets.Item(14)
$wsDiagnostics.Activate() | Out-Null
